# "Running all the suites by making runmode as Y in the excel sheets"
#
# On the "Test Cases" sheet, force every test case's Runmode to "Y" so the
# whole suite runs, and reset its Results column back to "SKIP" (no result
# has been recorded yet, since the suite hasn't actually run).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Test Cases")

$lastRow = $ws.UsedRange.Rows.Count   # header row + data rows

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells($r, 3).Value = "Y"       # column C: Runmode
    $ws.Cells($r, 4).Value = "SKIP"    # column D: Results
}
